# News From The Forest Links.xlsx — add the "January 2020" newsletter row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row under the existing list (row 16 = "November 2019")
$ws.Range("A17").Value = "January 2020"
$ws.Range("B17").Value = "https://myemail.constantcontact.com/News-From-The-Forest---January.html?soid=1102494320279&aid=W-tzxnGAuTU"

# Hyperlink the new URL cell, then restore the standard Hyperlink style
# (Hyperlinks.Add otherwise mints a fresh style index)
$ws.Hyperlinks.Add($ws.Range("B17"), "https://myemail.constantcontact.com/News-From-The-Forest---January.html?soid=1102494320279&aid=W-tzxnGAuTU")
$ws.Range("B17").Style = "Hyperlink"

# The author's selection ended up on B23 after adding the new row
$ws.Range("B23").Select() | Out-Null
